# issue #5: stock data output to json file
# Insert a new "property_category" column (with value "stock") into the
# 股票 (stock) worksheet, between the "total" and "date" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Shift the existing "date", "legislator_name", "legislator_id" columns
# (H, I, J) one column to the right (I, J, K) to make room for the new
# "property_category" column. Work right-to-left so a source column is
# always copied before it gets overwritten.
$ws.Range("J1:J2").Copy($ws.Range("K1"))
$ws.Range("I1:I2").Copy($ws.Range("J1"))
$ws.Range("H1:H2").Copy($ws.Range("I1"))

# New header cell (reuse the formatting of the neighboring "total" header)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "property_category"

# New data cell (reuse the formatting of the neighboring "total" value)
$ws.Range("G2").Copy($ws.Range("H2"))
$ws.Range("H2").Value = "stock"
